# Data-driven testing update: swap in a fresh set of login credentials on
# DATA_SHEET, trim the now-stale result rows from OUTPUT_DATASHEET, and
# move the active selection to reflect where the tester is currently working.

$wb = $excel.ActiveWorkbook

# --- DATA_SHEET: update the second data row's email + password values ---
$ws1 = $wb.Worksheets.Item("DATA_SHEET")
$ws1.Range("C2").Value = "DUVAKESH123@GMAIL.COM"
$ws1.Range("D2").Value = "DUVA999888"
$ws1.Range("E2").Value = "DUVA999888"

# --- OUTPUT_DATASHEET: clear out the stale sample result rows, keep header ---
$ws2 = $wb.Worksheets.Item("OUTPUT_DATASHEET")
$ws2.Rows("2:5").Delete()

# --- Update active selections/tab to match the current working state ---
$ws2.Activate()
$ws2.Range("A2:C5").Select()

$ws1.Activate()
$ws1.Range("F2").Select()
